# LOGGIN.xlsx - "solucion de interconexion que presentaba problemas"
#
# Appends 20 new login/logout session rows (41-60) to the "Logins" sheet,
# all belonging to role id "1" / "ADMINISTRADOR". Mirrors the existing
# rows: column A/B hold "dd/MM/yyyy HH:mm:ss" timestamps (stored as text,
# same visual style as the current data), columns C/D hold the role id and
# role name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the target ranges so the new cells line up with the existing
# columns: A/B keep the "Tiempo Entrada/Salida" date-time display, C is
# forced to text so the numeric-looking role id "1" isn't turned into a
# number by the Value assignment below (same as the existing "2"/"3" role
# ids elsewhere in the sheet, which are stored as text).
$ws.Range("A41:B60").NumberFormat = "dd/MM/yyyy HH:mm:ss"
$ws.Range("C41:C60").NumberFormat = "@"

# r, Tiempo Entrada, Tiempo Salida
$sessions = @(
    @(41, "27/10/2024 10:04:28", "27/10/2024 10:10:36"),
    @(42, "27/10/2024 08:28:38", "27/10/2024 08:28:53"),
    @(43, "27/10/2024 08:29:51", "27/10/2024 08:30:06"),
    @(44, "27/10/2024 08:44:52", "27/10/2024 08:45:34"),
    @(45, "27/10/2024 08:46:27", "27/10/2024 08:47:56"),
    @(46, "27/10/2024 08:49:30", "27/10/2024 08:50:59"),
    @(47, "27/10/2024 08:51:48", "27/10/2024 08:52:07"),
    @(48, "27/10/2024 08:52:57", "27/10/2024 08:55:57"),
    @(49, "27/10/2024 08:59:01", "27/10/2024 08:59:19"),
    @(50, "27/10/2024 08:59:54", "27/10/2024 09:01:43"),
    @(51, "27/10/2024 09:02:51", "27/10/2024 09:03:08"),
    @(52, "27/10/2024 09:03:52", "27/10/2024 09:06:44"),
    @(53, "27/10/2024 09:13:19", "27/10/2024 09:13:58"),
    @(54, "27/10/2024 09:14:35", "27/10/2024 09:18:29"),
    @(55, "27/10/2024 09:19:24", "27/10/2024 09:21:12"),
    @(56, "27/10/2024 09:28:41", "27/10/2024 09:31:49"),
    @(57, "27/10/2024 09:35:22", "27/10/2024 09:37:30"),
    @(58, "27/10/2024 09:41:47", "27/10/2024 09:45:29"),
    @(59, "27/10/2024 09:46:11", "27/10/2024 09:46:55"),
    @(60, "27/10/2024 10:04:28", "27/10/2024 10:10:36")
)

foreach ($session in $sessions) {
    $row = $session[0]
    $ws.Range("A$row").Value = $session[1]
    $ws.Range("B$row").Value = $session[2]
    $ws.Range("C$row").Value = "1"
    $ws.Range("D$row").Value = "ADMINISTRADOR"
}
